$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1, matching style of existing header cells (E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Per-row time_taken values (plain style, matching E2:E9 cells)
$timestamps = @(
    "2021-10-05 10:50:27.818509",
    "2021-10-05 10:50:27.818522",
    "2021-10-05 10:50:27.818526",
    "2021-10-05 10:50:27.818530",
    "2021-10-05 10:50:27.818533",
    "2021-10-05 10:50:27.818536",
    "2021-10-05 10:50:27.818570",
    "2021-10-05 10:50:27.818591"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
